# Revert the multi-industry "AI/ML" template branding back to the
# generic IT / Cloud Infrastructure wording on the "Change Management
# Overview" and "Change Impact Assessment" sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Change Management Overview")
$wsImpact   = $wb.Worksheets.Item("Change Impact Assessment")

# --- Change Management Overview ---------------------------------------
$wsOverview.Range("B6").Value = "Enterprise Cloud Infrastructure Migration"
$wsOverview.Range("A15").Value = "1. Achieve 95% user adoption of new IT systems within 6 months of go-live"
$wsOverview.Range("A17").Value = "3. Build organizational capability and confidence in IT technologies"
$wsOverview.Range("A20").Value = "6. Create positive stakeholder sentiment and enthusiasm for IT transformation"

# --- Change Impact Assessment ------------------------------------------
$wsImpact.Range("A4").Value = "IT Managers"
$wsImpact.Range("G4").Value = "IT automation"
$wsImpact.Range("A5").Value = "System Administrators"
